# Edit script for LOM3260.xlsx
# Applies the shared-strings/content restructuring described by the commit diff:
#  - adds PT objectives / PT short+full syllabus / bibliography text
#  - inserts a dedicated "Docentes responsaveis" pair of name rows (13/14)
#  - shifts label/content rows 15-21 down to accommodate, and appends new
#    rows 22 (Norma de recuperacao) and 23 (Bibliografia)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- literal text for the cells that actually change ---
$s18 = @'
Fornecer ao aluno uma introdução à computação científica moderna, usando a linguagem Python e suas bibliotecas numéricas e gráficas mais populares: numpy, scipy, matplotlib e pandas. Ao final do curso, o aluno estará capacitado a desenvolver programas complexos, de pequeno e médio porte para solucionar problemas de engenharia que envolvam processamento numérico de grandes conjuntos de dados e correlacionar variáveis usando métodos numéricos.
'@
$s22 = @'
7290967 - Emerson Gonçalves de Melo
'@
$s23 = @'
1176388 - Luiz Tadeu Fernandes Eleno
'@
$s24 = @'
Programa resumido:
'@
$s25 = @'
Introdução à programação em Python; palavras-chave em Python; rotinas e funções; classes; numpy e o conceito de slicing e indexing de arrays; revisão de métodos numéricos usando scipy; geração de gráficos e animações com a biblioteca matplotlib; criação de interfaces gráficas com o usuário usando matplotlib.widgets
'@
$s26 = @'
Short syllabus:
'@
$s27 = @'
Introduction to Python programming; keywords in Python; routines and functions; classes; numpy and the concept of slicing and indexing arrays; review of numerical methods using scipy; generating graphics and animations with the matplotlib library; creating graphical user interfaces using matplotlib.widgets
'@
$s28 = @'
Programa:
'@
$s29 = @'
• Introdução à programação em Python • Instalação de uma distribuição Python em Windows e Linux • Formatação de arquivos em Python • Estruturas condicionais • Laços de repetição de comandos • Outras palavras-chaves e métodos • Rotinas e funções • Códigos multifonte e bibliotecas pessoais • Bibliotecas numéricas e gráficas: numpy, scipy e matplotlib • Programação orientada a objeto: classes • Conceito de objetos e instâncias • Classes e subclasses• “Arrays” em numpy • O conceito de array em numpy • “Slicing” e indexação • Trabalhando com arquivos (entrada e saída) • Gráficos em matplotlib • A biblioteca matplotlib.pyplot e gráficos em 2D e 3D • A biblioteca matplotlib.animation para criar gráficos animados. • Interfaces gráficas com o usuário (Graphical User Interface, GUI) • Interfaces simples com a biblioteca matplotlib.widgets.
'@
$s30 = @'
Syllabus:
'@
$s31 = @'
• Introduction to Python programming • Installing a Python distribution on Windows and Linux • Python file formatting • Conditional structures • Command loops • Other keywords and methods • Routines and functions • Multi-source code and personal libraries • Numerical and graphical libraries: numpy, scipy and matplotlib • Object-oriented programming: classes • Concept of objects and instances • Classes and subclasses • “Arrays” in numpy • The concept of array in numpy • “Slicing” and indexing • Working with files ( input and output) • Graphs in matplotlib • The matplotlib.pyplot library and 2D and 3D graphs • The matplotlib.animation library for creating animated graphs. • Graphical User Interface (GUI) • Simple interfaces with the matplotlib.widgets library.
'@
$s32 = @'
Avaliação:
'@
$s33 = @'
Método:
'@
$s34 = @'
Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados.
'@
$s35 = @'
Critério:
'@
$s36 = @'
Média aritmética de exercícios e trabalhos propostos ao longo do curso e uma apresentação final de projeto.
'@
$s37 = @'
Norma de recuperação:
'@
$s38 = @'
Não haverá exame de recuperação.
'@
$s39 = @'
Bibliografia:
'@
$s40 = @'
Lambert, K. A. Fundamentos de Python: estruturas de dados. Cengage, 2ed, 2022.Nilo Ney Coutinho Menezes. Introdução à Programação com Python: Algoritmos e Lógica de Programação Para Iniciantes, 3a ed, 2019.Ramalho, L. Python Fluente. O’Reilly-Novatec, 2015Downey, A. B. Pense em Python. O’Reilly-Novatec, 2016.STEWART, J. M. Python for scientists. Cambridge University Press, 2014.TELLES, M. Python Power, Boston: Thomson Course Technology PTR, 2008.LUTZ, Mark. Programming Python, 3a ed, Sebastopol, CA: O’Reilly Media, 2006.MCGREGGOR, D. M. Mastering matplotlib. Birmingham, UK: Packt Publishing, 2015.
'@

# --- row 10: objectives content (PT/EN) replaces the old Emerson-name placeholder ---
$ws.Cells.Item(10, 2).Value = $s18
$ws.Cells.Item(10, 3).Value = $s18

# --- row 12 "Docentes responsaveis:" label already correct; add name rows 13/14 ---
$ws.Cells.Item(13, 1).ClearContents()
$ws.Cells.Item(13, 2).Value = $s22
$ws.Cells.Item(13, 3).Value = $s22
$ws.Cells.Item(14, 1).ClearContents()
$ws.Cells.Item(14, 2).Value = $s23
$ws.Cells.Item(14, 3).Value = $s23

# --- rows 15-18: resumo/short-syllabus/programa/syllabus, each now on its own row ---
$ws.Cells.Item(15, 1).Value = $s24
$ws.Cells.Item(15, 2).Value = $s25
$ws.Cells.Item(15, 3).Value = $s25
$ws.Cells.Item(16, 1).Value = $s26
$ws.Cells.Item(16, 2).Value = $s27
$ws.Cells.Item(16, 3).Value = $s27
$ws.Cells.Item(17, 1).Value = $s28
$ws.Cells.Item(17, 2).Value = $s29
$ws.Cells.Item(17, 3).Value = $s29
$ws.Cells.Item(18, 1).Value = $s30
$ws.Cells.Item(18, 2).Value = $s31
$ws.Cells.Item(18, 3).Value = $s31

# --- row 19: "Avaliacao:" label only (clear old Aulas-expositivas content) ---
$ws.Cells.Item(19, 1).Value = $s32
$ws.Cells.Item(19, 2).ClearContents()
$ws.Cells.Item(19, 3).ClearContents()

# --- rows 20-21: metodo / criterio ---
$ws.Cells.Item(20, 1).Value = $s33
$ws.Cells.Item(20, 2).Value = $s34
$ws.Cells.Item(20, 3).Value = $s34
$ws.Cells.Item(21, 1).Value = $s35
$ws.Cells.Item(21, 2).Value = $s36
$ws.Cells.Item(21, 3).Value = $s36

# --- new rows 22-23: norma de recuperacao / bibliografia ---
$ws.Cells.Item(22, 1).Value = $s37
$ws.Cells.Item(22, 2).Value = $s38
$ws.Cells.Item(22, 3).Value = $s38
$ws.Cells.Item(23, 1).Value = $s39
$ws.Cells.Item(23, 2).Value = $s40
$ws.Cells.Item(23, 3).Value = $s40

# --- newly-created rows 22/23 inherit a stray style in column B; re-apply the
#     correct (wrap-text, non-bold) format copied from the column B template cell ---
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row heights: shrink 15/16/21 to 60, grow 17/18 to 120, clear custom height
#     on 13/14/19, and set heights for the two appended rows ---
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(19).AutoFit()

